# "Plano de ação att" - update the Plano de Ação action-plan sheet:
#  - bump SPRINT progress % on a couple of rows
#  - fill in owner/priority on two previously-blank action rows
#  - clear out the leftover sample/template row (rows 49-51)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plano de Ação")
$ws.Activate()

# Row 39 ("Dashboard Estática"): progress 60% -> 95%
$ws.Range("E39").Value = 0.95

# Row 42 ("tela cadastrar maquina"): progress left blank -> 0%
$ws.Range("E42").NumberFormat = "0%"
$ws.Range("E42").Value = 0

# Row 43 ("BD VM" / "Tela Redefinir senha"): progress left blank -> 0%
$ws.Range("E43").NumberFormat = "0%"
$ws.Range("E43").Value = 0

# Row 44: assign responsible + progress
$ws.Range("C44").Value = "Gabriela "
$ws.Range("E44").NumberFormat = "0%"
$ws.Range("E44").Value = 0

# Row 45: assign responsible + priority + progress
$ws.Range("C45").Value = "Gabriela "
$ws.Range("D45").Value = "Essencial"
$ws.Range("E45").NumberFormat = "0%"
$ws.Range("E45").Value = 0.3

# Rows 49-51 were leftover sample data ("Todos da equipe" / "Criação da
# dashboard..." / etc.) - clear it back out, keeping the formatting.
$ws.Range("C49:H51").ClearContents()

# Leave the selection where the user left it after clearing that block.
$ws.Range("C49:H51").Select()
